$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")
$ws.Activate()
$ws.Range("H35").Value = 12
$ws.Range("H35").Select() | Out-Null
